$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray row (row 5, which only held a lone " " value in column G)
$ws.Rows.Item(5).Delete()

# Rename the worksheet tab from "Hoja1" to "Importar"
$ws.Name = "Importar"

# Move the cell selection/active cell to B10 (was K3)
$ws.Range("B10").Select()
